$d = $word.ActiveDocument
$totalFound = 0

$old0 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:870)'
$new0 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1074)'
$found0 = $d.Content.Find.Execute($old0, $true, $false, $false, $false, $false, $true, 1, $false, $new0, 1)
if (-not $found0) { Write-Output "WARNING: replacement 0 not found" }
Write-Output "Replace 0: $found0"

$old1 = '	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:150)'
$new1 = '	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:160)'
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 1)
if (-not $found1) { Write-Output "WARNING: replacement 1 not found" }
Write-Output "Replace 1: $found1"

$old2 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1034)'
$new2 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)'
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 1)
if (-not $found2) { Write-Output "WARNING: replacement 2 not found" }
Write-Output "Replace 2: $found2"

$old3 = '	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:183)'
$new3 = '	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)'
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 1)
if (-not $found3) { Write-Output "WARNING: replacement 3 not found" }
Write-Output "Replace 3: $found3"

$old4 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:297)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:1)
	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:201)'
$new4 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)
	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)'
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 1)
if (-not $found4) { Write-Output "WARNING: replacement 4 not found" }
Write-Output "Replace 4: $found4"

$old5 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:259)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)
	at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:246)
	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)
	at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:836)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:252)
	at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:691)
	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:396)
	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:318)
	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)
	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)'
$new5 = '	at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)
	at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)
	at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)
	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)
	at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)
	at sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)'
$found5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 1)
if (-not $found5) { Write-Output "WARNING: replacement 5 not found" }
Write-Output "Replace 5: $found5"

$old6 = '	at org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)'
$new6 = '	at org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)
	at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)'
$found6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 1)
if (-not $found6) { Write-Output "WARNING: replacement 6 not found" }
Write-Output "Replace 6: $found6"

$old7 = '	at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)
'
$new7 = ''
$found7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 1)
if (-not $found7) { Write-Output "WARNING: replacement 7 not found" }
Write-Output "Replace 7: $found7"

$old8 = '	at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)
	at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)'
$new8 = '	at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)
	at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)
	at org.junit.runners.ParentRunner.run(ParentRunner.java:363)
	at org.junit.runners.Suite.runChild(Suite.java:128)
	at org.junit.runners.Suite.runChild(Suite.java:27)
	at org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)
	at org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)
	at org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)
	at org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)
	at org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)'
$found8 = $d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 1)
if (-not $found8) { Write-Output "WARNING: replacement 8 not found" }
Write-Output "Replace 8: $found8"

$old9 = '	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:459)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:675)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:382)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:192)'
$new9 = '	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)
	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)'
$found9 = $d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 1)
if (-not $found9) { Write-Output "WARNING: replacement 9 not found" }
Write-Output "Replace 9: $found9"

